{"js": "// Update the worksheet date and all 25 \"two-digit \u00d7 two-digit\" answer\n// cells to the new day's values. Every \"before\" string below is unique\n// within the document, so a plain search + Replace on each is safe and\n// keeps all run/paragraph formatting (font, size, etc.) untouched.\nconst replacements = [\n  [\"2024-04-21 Sunday\", \"2024-04-22 Monday\"],\n  [\"85\u00d768=5780\", \"77\u00d776=5852\"],\n  [\"36\u00d750=1800\", \"24\u00d717=408\"],\n  [\"64\u00d720=1280\", \"49\u00d787=4263\"],\n  [\"14\u00d779=1106\", \"37\u00d727=999\"],\n  [\"44\u00d779=3476\", \"71\u00d755=3905\"],\n  [\"88\u00d782=7216\", \"35\u00d785=2975\"],\n  [\"71\u00d769=4899\", \"44\u00d799=4356\"],\n  [\"85\u00d728=2380\", \"28\u00d785=2380\"],\n  [\"96\u00d772=6912\", \"98\u00d734=3332\"],\n  [\"79\u00d795=7505\", \"54\u00d784=4536\"],\n  [\"47\u00d711=517\", \"32\u00d781=2592\"],\n  [\"52\u00d793=4836\", \"42\u00d772=3024\"],\n  [\"80\u00d789=7120\", \"72\u00d740=2880\"],\n  [\"76\u00d756=4256\", \"54\u00d769=3726\"],\n  [\"24\u00d719=456\", \"66\u00d725=1650\"],\n  [\"42\u00d720=840\", \"63\u00d737=2331\"],\n  [\"35\u00d794=3290\", \"42\u00d776=3192\"],\n  [\"41\u00d792=3772\", \"82\u00d764=5248\"],\n  [\"35\u00d755=1925\", \"73\u00d732=2336\"],\n  [\"72\u00d735=2520\", \"79\u00d742=3318\"],\n  [\"74\u00d714=1036\", \"52\u00d768=3536\"],\n  [\"48\u00d750=2400\", \"14\u00d767=938\"],\n  [\"98\u00d747=4606\", \"40\u00d781=3240\"],\n  [\"77\u00d764=4928\", \"55\u00d768=3740\"],\n  [\"34\u00d721=714\", \"82\u00d760=4920\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and all 25 \"two-digit x two-digit\" answer\n# cells to the new day's values. Every \"before\" string is unique within\n# the document, so Find/Replace (scoped to a single hit) on each pair is\n# safe and leaves all run/paragraph formatting (font, size, etc.) intact.\n\n$replacements = @(\n    @(\"2024-04-21 Sunday\", \"2024-04-22 Monday\"),\n    @(\"85\u00d768=5780\", \"77\u00d776=5852\"),\n    @(\"36\u00d750=1800\", \"24\u00d717=408\"),\n    @(\"64\u00d720=1280\", \"49\u00d787=4263\"),\n    @(\"14\u00d779=1106\", \"37\u00d727=999\"),\n    @(\"44\u00d779=3476\", \"71\u00d755=3905\"),\n    @(\"88\u00d782=7216\", \"35\u00d785=2975\"),\n    @(\"71\u00d769=4899\", \"44\u00d799=4356\"),\n    @(\"85\u00d728=2380\", \"28\u00d785=2380\"),\n    @(\"96\u00d772=6912\", \"98\u00d734=3332\"),\n    @(\"79\u00d795=7505\", \"54\u00d784=4536\"),\n    @(\"47\u00d711=517\", \"32\u00d781=2592\"),\n    @(\"52\u00d793=4836\", \"42\u00d772=3024\"),\n    @(\"80\u00d789=7120\", \"72\u00d740=2880\"),\n    @(\"76\u00d756=4256\", \"54\u00d769=3726\"),\n    @(\"24\u00d719=456\", \"66\u00d725=1650\"),\n    @(\"42\u00d720=840\", \"63\u00d737=2331\"),\n    @(\"35\u00d794=3290\", \"42\u00d776=3192\"),\n    @(\"41\u00d792=3772\", \"82\u00d764=5248\"),\n    @(\"35\u00d755=1925\", \"73\u00d732=2336\"),\n    @(\"72\u00d735=2520\", \"79\u00d742=3318\"),\n    @(\"74\u00d714=1036\", \"52\u00d768=3536\"),\n    @(\"48\u00d750=2400\", \"14\u00d767=938\"),\n    @(\"98\u00d747=4606\", \"40\u00d781=3240\"),\n    @(\"77\u00d764=4928\", \"55\u00d768=3740\"),\n    @(\"34\u00d721=714\", \"82\u00d760=4920\")\n)\n\n$d = $word.ActiveDocument\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n\n    # MatchCase:$true, Forward:$true, Wrap:wdFindContinue(1), Format:$false,\n    # Replace:wdReplaceAll(2) -- each string occurs exactly once so this is\n    # equivalent to a single targeted replace.\n    $find.Execute(\n        $oldText,   # FindText\n        $true,      # MatchCase\n        $false,     # MatchWholeWord\n        $false,     # MatchWildcards\n        $false,     # MatchSoundsLike\n        $false,     # MatchAllWordForms\n        $true,      # Forward\n        1,          # Wrap (wdFindContinue)\n        $false,     # Format\n        $newText,   # ReplaceWith\n        2           # Replace (wdReplaceAll)\n    ) | Out-Null\n}\n"}
